# Updates horarios data for Linea 141 (scrape refresh 08:38:27 -> 08:52:20)
# - Refreshes timestamps / row counts
# - Corrects several Hora_Scrap / Linea / Minutos values across the three sheets
# - Appends newly scraped rows 119-126 on sheet "LP1912"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet: LP1912 ---
$ws1.Cells.Item(2, 1).Value = 'Última actualización: 08:52:20'
$ws1.Cells.Item(3, 1).Value = 'Total filas: 121'
$ws1.Cells.Item(43, 1).Value = '06:24:16'
$ws1.Cells.Item(43, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(43, 4).Value = 52
$ws1.Cells.Item(44, 1).Value = '06:53:31'
$ws1.Cells.Item(44, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(44, 4).Value = 23
$ws1.Cells.Item(47, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(48, 3).Value = '10_OLMOS'
$ws1.Cells.Item(52, 1).Value = '05:54:55'
$ws1.Cells.Item(52, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(52, 4).Value = 97
$ws1.Cells.Item(53, 1).Value = '06:53:31'
$ws1.Cells.Item(53, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(53, 4).Value = 38
$ws1.Cells.Item(54, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(56, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(65, 1).Value = '07:50:27'
$ws1.Cells.Item(65, 3).Value = '10_OLMOS'
$ws1.Cells.Item(65, 4).Value = 2
$ws1.Cells.Item(66, 1).Value = '07:18:07'
$ws1.Cells.Item(66, 3).Value = '215D_EL PATO'
$ws1.Cells.Item(66, 4).Value = 34
$ws1.Cells.Item(78, 1).Value = '07:18:07'
$ws1.Cells.Item(78, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(78, 4).Value = 65
$ws1.Cells.Item(79, 1).Value = '08:16:28'
$ws1.Cells.Item(79, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(79, 4).Value = 7
$ws1.Cells.Item(89, 1).Value = '08:52:20'
$ws1.Cells.Item(89, 2).Value = '08:52'
$ws1.Cells.Item(89, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(89, 4).Value = 0
$ws1.Cells.Item(90, 1).Value = '08:52:20'
$ws1.Cells.Item(90, 2).Value = '08:52'
$ws1.Cells.Item(90, 3).Value = '10_OLMOS'
$ws1.Cells.Item(90, 4).Value = 0
$ws1.Cells.Item(91, 2).Value = '08:53'
$ws1.Cells.Item(91, 3).Value = '10_OLMOS'
$ws1.Cells.Item(91, 4).Value = 15
$ws1.Cells.Item(92, 1).Value = '08:52:20'
$ws1.Cells.Item(92, 2).Value = '08:54'
$ws1.Cells.Item(92, 3).Value = '17_ROMERO'
$ws1.Cells.Item(92, 4).Value = 2
$ws1.Cells.Item(93, 1).Value = '08:52:20'
$ws1.Cells.Item(93, 2).Value = '09:01'
$ws1.Cells.Item(93, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(93, 4).Value = 9
$ws1.Cells.Item(94, 1).Value = '07:18:07'
$ws1.Cells.Item(94, 2).Value = '09:02'
$ws1.Cells.Item(94, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(94, 4).Value = 104
$ws1.Cells.Item(95, 1).Value = '08:52:20'
$ws1.Cells.Item(95, 2).Value = '09:03'
$ws1.Cells.Item(95, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(95, 4).Value = 11
$ws1.Cells.Item(96, 1).Value = '08:16:28'
$ws1.Cells.Item(96, 2).Value = '09:04'
$ws1.Cells.Item(96, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(96, 4).Value = 48
$ws1.Cells.Item(97, 1).Value = '08:52:20'
$ws1.Cells.Item(97, 2).Value = '09:07'
$ws1.Cells.Item(97, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(97, 4).Value = 15
$ws1.Cells.Item(98, 1).Value = '08:52:20'
$ws1.Cells.Item(98, 2).Value = '09:10'
$ws1.Cells.Item(98, 4).Value = 18
$ws1.Cells.Item(99, 1).Value = '07:50:27'
$ws1.Cells.Item(99, 2).Value = '09:11'
$ws1.Cells.Item(99, 3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(99, 4).Value = 81
$ws1.Cells.Item(100, 1).Value = '07:18:07'
$ws1.Cells.Item(100, 2).Value = '09:11'
$ws1.Cells.Item(100, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(100, 4).Value = 113
$ws1.Cells.Item(101, 1).Value = '08:52:20'
$ws1.Cells.Item(101, 2).Value = '09:13'
$ws1.Cells.Item(101, 3).Value = '10_OLMOS'
$ws1.Cells.Item(101, 4).Value = 21
$ws1.Cells.Item(102, 1).Value = '08:52:20'
$ws1.Cells.Item(102, 2).Value = '09:15'
$ws1.Cells.Item(102, 3).Value = '14_ABASTO'
$ws1.Cells.Item(102, 4).Value = 23
$ws1.Cells.Item(103, 2).Value = '09:16'
$ws1.Cells.Item(103, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(103, 4).Value = 38
$ws1.Cells.Item(104, 1).Value = '08:52:20'
$ws1.Cells.Item(104, 2).Value = '09:17'
$ws1.Cells.Item(104, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(104, 4).Value = 25
$ws1.Cells.Item(105, 1).Value = '08:52:20'
$ws1.Cells.Item(105, 2).Value = '09:21'
$ws1.Cells.Item(105, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(105, 4).Value = 29
$ws1.Cells.Item(106, 2).Value = '09:22'
$ws1.Cells.Item(106, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(106, 4).Value = 44
$ws1.Cells.Item(107, 2).Value = '09:22'
$ws1.Cells.Item(107, 3).Value = '17_ROMERO'
$ws1.Cells.Item(107, 4).Value = 66
$ws1.Cells.Item(108, 1).Value = '08:52:20'
$ws1.Cells.Item(108, 2).Value = '09:23'
$ws1.Cells.Item(108, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(108, 4).Value = 31
$ws1.Cells.Item(109, 2).Value = '09:23'
$ws1.Cells.Item(109, 3).Value = '17_ROMERO'
$ws1.Cells.Item(109, 4).Value = 45
$ws1.Cells.Item(110, 1).Value = '08:52:20'
$ws1.Cells.Item(110, 2).Value = '09:23'
$ws1.Cells.Item(110, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(110, 4).Value = 31
$ws1.Cells.Item(111, 2).Value = '09:29'
$ws1.Cells.Item(111, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(111, 4).Value = 73
$ws1.Cells.Item(112, 1).Value = '08:52:20'
$ws1.Cells.Item(112, 2).Value = '09:32'
$ws1.Cells.Item(112, 3).Value = '15_ABASTO'
$ws1.Cells.Item(112, 4).Value = 40
$ws1.Cells.Item(113, 1).Value = '08:52:20'
$ws1.Cells.Item(113, 2).Value = '09:33'
$ws1.Cells.Item(113, 3).Value = '10_OLMOS'
$ws1.Cells.Item(113, 4).Value = 41
$ws1.Cells.Item(114, 1).Value = '08:38:27'
$ws1.Cells.Item(114, 2).Value = '09:34'
$ws1.Cells.Item(114, 4).Value = 56
$ws1.Cells.Item(115, 1).Value = '08:52:20'
$ws1.Cells.Item(115, 2).Value = '09:35'
$ws1.Cells.Item(115, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(115, 4).Value = 43
$ws1.Cells.Item(116, 1).Value = '08:52:20'
$ws1.Cells.Item(116, 2).Value = '09:35'
$ws1.Cells.Item(116, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(116, 4).Value = 43
$ws1.Cells.Item(117, 1).Value = '08:16:28'
$ws1.Cells.Item(117, 2).Value = '09:37'
$ws1.Cells.Item(117, 3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(117, 4).Value = 81
$ws1.Cells.Item(118, 1).Value = '08:52:20'
$ws1.Cells.Item(118, 2).Value = '09:42'
$ws1.Cells.Item(118, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(118, 4).Value = 50
$ws1.Cells.Item(119, 1).Value = '08:52:20'
$ws1.Cells.Item(119, 2).Value = '09:43'
$ws1.Cells.Item(119, 3).Value = '14_ABASTO'
$ws1.Cells.Item(119, 4).Value = 51
$ws1.Cells.Item(119, 5).Value = 'LP1912'
$ws1.Cells.Item(120, 1).Value = '08:16:28'
$ws1.Cells.Item(120, 2).Value = '09:46'
$ws1.Cells.Item(120, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(120, 4).Value = 90
$ws1.Cells.Item(120, 5).Value = 'LP1912'
$ws1.Cells.Item(121, 1).Value = '08:52:20'
$ws1.Cells.Item(121, 2).Value = '09:52'
$ws1.Cells.Item(121, 3).Value = '15_ABASTO'
$ws1.Cells.Item(121, 4).Value = 60
$ws1.Cells.Item(121, 5).Value = 'LP1912'
$ws1.Cells.Item(122, 1).Value = '08:52:20'
$ws1.Cells.Item(122, 2).Value = '10:10'
$ws1.Cells.Item(122, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(122, 4).Value = 78
$ws1.Cells.Item(122, 5).Value = 'LP1912'
$ws1.Cells.Item(123, 1).Value = '08:38:27'
$ws1.Cells.Item(123, 2).Value = '10:12'
$ws1.Cells.Item(123, 3).Value = '15_ABASTO'
$ws1.Cells.Item(123, 4).Value = 94
$ws1.Cells.Item(123, 5).Value = 'LP1912'
$ws1.Cells.Item(124, 1).Value = '08:52:20'
$ws1.Cells.Item(124, 2).Value = '10:21'
$ws1.Cells.Item(124, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(124, 4).Value = 89
$ws1.Cells.Item(124, 5).Value = 'LP1912'
$ws1.Cells.Item(125, 1).Value = '08:52:20'
$ws1.Cells.Item(125, 2).Value = '10:26'
$ws1.Cells.Item(125, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(125, 4).Value = 94
$ws1.Cells.Item(125, 5).Value = 'LP1912'
$ws1.Cells.Item(126, 1).Value = '08:52:20'
$ws1.Cells.Item(126, 2).Value = '10:42'
$ws1.Cells.Item(126, 3).Value = '17_ROMERO'
$ws1.Cells.Item(126, 4).Value = 110
$ws1.Cells.Item(126, 5).Value = 'LP1912'

# --- Sheet: LP1912-215 ---
$ws2.Cells.Item(2, 1).Value = 'Última actualización: 08:52:20'
$ws2.Cells.Item(15, 1).Value = '08:52:20'
$ws2.Cells.Item(15, 4).Value = 9
$ws2.Cells.Item(17, 1).Value = '08:52:20'
$ws2.Cells.Item(17, 4).Value = 50
$ws2.Cells.Item(18, 1).Value = '08:52:20'
$ws2.Cells.Item(18, 4).Value = 94

# --- Sheet: 6203-6173 ---
$ws3.Cells.Item(2, 1).Value = 'Última actualización: 08:52:20'
$ws3.Cells.Item(21, 1).Value = '08:52:20'
$ws3.Cells.Item(21, 4).Value = 17
$ws3.Cells.Item(22, 1).Value = '08:52:20'
$ws3.Cells.Item(22, 4).Value = 71
